# Update "Horarios" workbook with the latest scrape (02:49:42) for Línea 141.
$wb = $excel.ActiveWorkbook

$oldTime = "02:37:48"
$newTime = "02:49:42"

# ---------------------------------------------------------------------------
# Sheet "LP1912": 3 existing data rows get their Hora_Scrap / Minutos
# refreshed, and a new row (215A_EL PATO) is appended.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Range("A6").Value = $newTime
$ws1.Range("D6").Value = 9

$ws1.Range("A7").Value = $newTime
$ws1.Range("D7").Value = 59

$ws1.Range("A8").Value = $newTime
$ws1.Range("D8").Value = 72

$ws1.Range("A9").Value = $newTime
$ws1.Range("B9").Value = "04:45"
$ws1.Range("C9").Value = "215A_EL PATO"
$ws1.Range("D9").Value = 116
$ws1.Range("E9").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": 1 existing data row refreshed, plus the same new
# 215A_EL PATO row appended.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Range("A6").Value = $newTime
$ws2.Range("D6").Value = 9

$ws2.Range("A7").Value = $newTime
$ws2.Range("B7").Value = "04:45"
$ws2.Range("C7").Value = "215A_EL PATO"
$ws2.Range("D7").Value = 116
$ws2.Range("E7").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": no schedule rows, only the "last updated" stamp moves.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
